# Swap the "Valor Mora" values between row 16 and row 24 in column F
# (F16: 30285 -> 21805, F24: 21805 -> 30285)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F16").Value = 21805
$ws.Range("F24").Value = 30285
